$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.226.11"
$ws.Range("E2").Value = "  -0.29%  "
$ws.Range("D3").Value = "3.405.73"
$ws.Range("E3").Value = "  +0.78%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.07"
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.09"
$ws.Range("E6").Value = "  +1.44%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.404.72"
$ws.Range("E8").Value = "  +0.83%  "
$ws.Range("E9").Value = "  -0.73%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.75"
$ws.Range("E10").Value = "  +3.70%  "
$ws.Range("E11").Value = "  -1.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.380"
$ws.Range("E12").Value = "  -1.98%  "
$ws.Range("D13").Value = "3.987.45"
$ws.Range("E13").Value = "  +0.88%  "
$ws.Range("E14").Value = "  -0.99%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.62"
$ws.Range("E15").Value = "  +1.79%  "
$ws.Range("E16").Value = "  -1.54%  "
$ws.Range("D17").Value = "3.403.62"
$ws.Range("E17").Value = "  +0.75%  "
$ws.Range("D18").Value = "61.324.39"
$ws.Range("E18").Value = "  -0.26%  "
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("E20").Value = "  -0.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.34"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "380.56"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("D23").Value = "3.535.35"
$ws.Range("E23").Value = "  +0.64%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("E25").Value = "  +0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.50"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000122"
$ws.Range("E27").Value = "  -2.09%  "
$ws.Range("E28").Value = "  +9.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.60"
$ws.Range("E29").Value = "  -8.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.45"
$ws.Range("E30").Value = "  -1.00%  "
$ws.Range("E31").Value = "  +0.13%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.14"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.46"
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("E36").Value = "  +1.49%  "
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("E38").Value = "  -1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "165.89"
$ws.Range("E39").Value = "  +0.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0774"
$ws.Range("E40").Value = "  +0.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.37"
$ws.Range("E41").Value = "  +8.13%  "
$ws.Range("E42").Value = "  +1.16%  "
$ws.Range("E43").Value = "  +0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.779"
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("E46").Value = "  +0.18%  "
$ws.Range("E47").Value = "  -2.94%  "
$ws.Range("D48").Value = "2.540.74"
$ws.Range("E48").Value = "  +7.40%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.87"
$ws.Range("E49").Value = "  +5.40%  "
$ws.Range("E50").Value = "  -0.60%  "
$ws.Range("E51").Value = "  -0.52%  "
